$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text, mirroring the source inlineStr cells
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "90.733.34"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "3.146.42"
$ws.Range("E3").Value = "  +1.77%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "217.31"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "626.37"
$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "1.12"
$ws.Range("E7").Value = "  +26.55%  "

$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "0.367"
$ws.Range("E8").Value = "  -3.06%  "

$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").Value = "3.145.37"
$ws.Range("E10").Value = "  +1.90%  "

$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "0.752"
$ws.Range("E11").Value = "  +5.97%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.200"
$ws.Range("E12").Value = "  +4.29%  "

$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  -2.18%  "

$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").Value = "5.71"
$ws.Range("E14").Value = "  +5.97%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "35.29"
$ws.Range("E15").Value = "  +6.83%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "90.557.88"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "3.745.02"
$ws.Range("E17").Value = "  +2.56%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.154.63"
$ws.Range("E18").Value = "  +2.41%  "

$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D19").Value = "3.73"
$ws.Range("E19").Value = "  +1.80%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "14.83"
$ws.Range("E20").Value = "  +7.20%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "478.06"
$ws.Range("E21").Value = "  +10.65%  "

$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").Value = "0.0000213"
$ws.Range("E22").Value = "  -7.20%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "9.32"
$ws.Range("E23").Value = "  +8.96%  "

$ws.Range("B24").Value = "Polkadot"
$ws.Range("C24").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D24").Value = "5.32"
$ws.Range("E24").Value = "  +3.17%  "

$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "5.95"
$ws.Range("E25").Value = "  +9.35%  "

$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "95.43"
$ws.Range("E26").Value = "  +13.81%  "

$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "12.38"
$ws.Range("E27").Value = "  +4.40%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "3.393.99"
$ws.Range("E28").Value = "  +4.57%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "9.45"
$ws.Range("E30").Value = "  +9.06%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "29.43"
$ws.Range("E31").Value = "  +26.74%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.217"
$ws.Range("E32").Value = "  +49.70%  "

$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").Value = "0.161"
$ws.Range("E33").Value = "  -5.49%  "

$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "513.78"
$ws.Range("E34").Value = "  -1.09%  "

$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").Value = "1.95"
$ws.Range("E35").Value = "  +5.98%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.147"
$ws.Range("E36").Value = "  +7.53%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "7.02"
$ws.Range("E37").Value = "  +0.98%  "

$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.59"
$ws.Range("E38").Value = "  -6.92%  "

$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "1.31"
$ws.Range("E39").Value = "  +2.02%  "

$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "0.435"
$ws.Range("E40").Value = "  +16.49%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "0.0867"
$ws.Range("E41").Value = "  +16.72%  "

$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").Value = "22.24"
$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("B43").Value = "Binance-PegBSC-USD"
$ws.Range("C43").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D43").Value = "0.760"
$ws.Range("E43").Value = "  +21.25%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "1.97"
$ws.Range("E45").Value = "  +5.72%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "4.78"
$ws.Range("E46").Value = "  +14.59%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "0.717"
$ws.Range("E47").Value = "  +18.12%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "151.29"
$ws.Range("E48").Value = "  +4.80%  "

$ws.Range("B49").Value = "MantraDAO"
$ws.Range("C49").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D49").Value = "3.02"
$ws.Range("E49").Value = "  +24.06%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "45.40"
$ws.Range("E50").Value = "  +3.45%  "

$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").Value = "1.35"
$ws.Range("E51").Value = "  +5.75%  "
